$d = $word.ActiveDocument

# --- Edit 1: update the date/time in the document's Date paragraph ---
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "March  16, 2022 (08:59:34 AM)`r") {
        $p.Range.Text = "March  16, 2022 (09:08:50 PM)"
        break
    }
}

# --- Edit 2: add a new bullet item after "only include code in text form..." ---
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "only include code in text form such that it can be copy-pasted for reuse`r") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $insertionPoint = $target.Range.Duplicate
    $insertionPoint.Collapse(0)

    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body>' +
           '<w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1010"/></w:numPr></w:pPr>' +
           '<w:r><w:t xml:space="preserve">only include code in text form such that it can be copy-pasted for reuse</w:t></w:r></w:p>' +
           '<w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1010"/></w:numPr></w:pPr>' +
           '<w:r><w:t xml:space="preserve">make sure to include blank lines before and after code blocks, since the absence of these can cause the code block to display incorrectly.</w:t></w:r></w:p>' +
           '</w:body></w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'

    $insertionPoint.InsertXML($xml)
}
